$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row
$ws.Range("A1").Value = "username"
$ws.Range("B1").Value = "password"

# Existing data rows remain the same (Raju/raju12, Bolu/bolu123)

# Add new rows
$ws.Range("A4").Value = "Dolu"
$ws.Range("B4").Value = "dolu123"
$ws.Range("A5").Value = "Rama"
$ws.Range("B5").Value = "rama123"

$ws.Range("B5").Select()
